# "apresentar tabela com validacao de perfil"
# - remove columns "gerar_etiqueta" and "data_carta_gerada"
# - rename "solicitar_correcao" -> "corrigido"
# - add new "retida" column
# - blank out the old "gerar_etiqueta" Sim/Nao flag (now just an empty date-formatted cell)
#   and fill the new "corrigido"/"retida" columns with "NAO" for every data row
# - scroll the sheet over and select P2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the two obsolete columns -----------------------------------
# column P = "gerar_etiqueta" (16th column)
$ws.Columns.Item(16).Delete()
# after the shift, the old "data_carta_gerada" column is now in position 17
$ws.Columns.Item(17).Delete()

# Columns are now:
#   P = numero_etiqueta, Q = quebra_sequencia, R = solicitar_correcao

# --- rename / add headers ------------------------------------------------
$ws.Range("R1").Value = "corrigido"
$ws.Range("S1").Value = "retida"

# --- data rows -------------------------------------------------------------
# P column keeps the date number-format used by the neighbouring columns
# (O = data_conferencia) but stays empty, just like before.
$ws.Range("O2").Copy()
$ws.Range("P2:P7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 18).Value = "NAO"   # corrigido
    $ws.Cells.Item($r, 19).Value = "NAO"   # retida
}

# --- cosmetics: column P width was manually widened (loses the bestFit autosize) ---
$ws.Columns.Item(16).ColumnWidth = 16.1667

# --- view state: scroll right and select P2 --------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("P2").Select()
